# Update imputed KNN values in the result data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B4"  = 7.593000000000001
    "B7"  = 5.024
    "C7"  = -13.383
    "C15" = -13.385
    "B16" = 4.826000000000001
    "E16" = 16.544
    "E19" = 16.593
    "C21" = -12.518
    "C22" = -12.961
    "C23" = -12.223
    "B28" = 5.709000000000001
    "B29" = 5.404000000000001
    "B32" = 6.388
    "C34" = -11.997
    "E36" = 16.486
    "B40" = 9.044999999999998
    "C43" = -13.516
    "C45" = -13.03
    "E46" = 16.934
    "C50" = -13.497
    "E50" = 16.542
    "C51" = -11.276
    "B52" = 5.503
    "B57" = 5.093000000000001
    "B66" = 5.054
    "C66" = -10.883
    "C67" = -11.38
    "C79" = -12.057
    "C84" = -13.72
    "C92" = -11.084
    "E95" = 17.399
    "C97" = -11.564
    "E97" = 17.114
    "B100" = 5.588
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
